$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: "Idioma" (Language) header, splitting existing data into
# "Japones" rows and duplicating them as "Mandarim" rows.
$ws.Range("E1").Value = "Idioma"

# Existing rows 2-11 are the Japanese entries.
$ws.Range("E2:E11").Value = "Japones"

# Duplicate rows 2-11 into rows 12-21 (same A:D content) for Mandarim.
for ($i = 2; $i -le 11; $i++) {
    $src = $i
    $dst = $i + 10
    $ws.Range("A$dst").Value = $ws.Range("A$src").Value2
    $ws.Range("B$dst").Value = $ws.Range("B$src").Value2
    $ws.Range("C$dst").Value = $ws.Range("C$src").Value2
    $ws.Range("D$dst").Value = $ws.Range("D$src").Value2
    $ws.Range("E$dst").Value = "Mandarim"
}

# Column E width.
$ws.Columns.Item(5).ColumnWidth = 9.67

# Empty, underlined cell at G17 (no value, just formatting) + selection there.
$ws.Range("G17").Font.Underline = $true
$ws.Range("G17").Select() | Out-Null

# Page setup.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
